# Auto-generated update script: applies numeric corrections to several
# leve-profit rows across multiple crafting-job sheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW), per the scheduled-runner data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 27: Don't Forget to Take Your Meds
$ws.Cells.Item(27, 8).Value = 0  # H27: 8000 -> 0
$ws.Cells.Item(27, 10).Value = 0  # J27: 8000 -> 0
$ws.Cells.Item(27, 12).Value = 0  # L27: 24000 -> 0
$ws.Cells.Item(27, 14).ClearContents()  # N27: -24202 -> (removed)

# Row 32: Automata for the People
$ws.Cells.Item(32, 8).Value = 3041.6667  # H32: 2942.5454 -> 3041.6667
$ws.Cells.Item(32, 9).Value = 3090.5  # I32: 2792.5 -> 3090.5
$ws.Cells.Item(32, 10).Value = 3027.7144  # J32: 3028.2856 -> 3027.7144
$ws.Cells.Item(32, 11).Value = 3090.5  # K32: 2792.5 -> 3090.5
$ws.Cells.Item(32, 12).Value = 3027.7144  # L32: 3028.2856 -> 3027.7144
$ws.Cells.Item(32, 13).Value = -2764.5  # M32: -2466.5 -> -2764.5
$ws.Cells.Item(32, 14).Value = -3679.7144  # N32: -3680.2856 -> -3679.7144

# Row 76: Warding Off Temptation
$ws.Cells.Item(76, 8).Value = 5000  # H76: 5333.3335 -> 5000
$ws.Cells.Item(76, 10).Value = 5000  # J76: 5333.3335 -> 5000
$ws.Cells.Item(76, 12).Value = 5000  # L76: 5333.3335 -> 5000
$ws.Cells.Item(76, 14).Value = -5630  # N76: -5963.3335 -> -5630

# Row 79: The Garden of Arcane Delights (L)
$ws.Cells.Item(79, 8).Value = 5000  # H79: 5333.3335 -> 5000
$ws.Cells.Item(79, 10).Value = 5000  # J79: 5333.3335 -> 5000
$ws.Cells.Item(79, 12).Value = 5000  # L79: 5333.3335 -> 5000
$ws.Cells.Item(79, 14).Value = -7184  # N79: -7517.3335 -> -7184

# Row 100: Asking for a Friend
$ws.Cells.Item(100, 8).Value = 2081.5  # H100: 1872.8889 -> 2081.5
$ws.Cells.Item(100, 9).Value = 1266.9231  # I100: 1229.6428 -> 1266.9231
$ws.Cells.Item(100, 10).Value = 4199.4  # J100: 4124.25 -> 4199.4
$ws.Cells.Item(100, 11).Value = 1266.9231  # K100: 1229.6428 -> 1266.9231
$ws.Cells.Item(100, 12).Value = 4199.4  # L100: 4124.25 -> 4199.4
$ws.Cells.Item(100, 13).Value = -725.9231  # M100: -688.6428000000001 -> -725.9231
$ws.Cells.Item(100, 14).Value = -5281.4  # N100: -5206.25 -> -5281.4

# Row 132: Fast-forwarding Flora
$ws.Cells.Item(132, 8).Value = 3945.3257  # H132: 4026.6428 -> 3945.3257
$ws.Cells.Item(132, 9).Value = 3284.3901  # I132: 3353.25 -> 3284.3901
$ws.Cells.Item(132, 11).Value = 9853.1703  # K132: 10059.75 -> 9853.1703
$ws.Cells.Item(132, 13).Value = -7323.1703  # M132: -7529.75 -> -7323.1703

# Row 137: Cutting Edge of Culinary Quality
$ws.Cells.Item(137, 8).Value = 6464.923  # H137: 5888.2144 -> 6464.923
$ws.Cells.Item(137, 9).Value = 1454.4  # I137: 1308.8572 -> 1454.4
$ws.Cells.Item(137, 10).Value = 9596.5  # J137: 10467.571 -> 9596.5
$ws.Cells.Item(137, 11).Value = 4363.200000000001  # K137: 3926.5716 -> 4363.200000000001
$ws.Cells.Item(137, 12).Value = 28789.5  # L137: 31402.713 -> 28789.5
$ws.Cells.Item(137, 13).Value = -1813.200000000001  # M137: -1376.5716 -> -1813.200000000001
$ws.Cells.Item(137, 14).Value = -33889.5  # N137: -36502.713 -> -33889.5

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Cells.Item(2, 8).Value = 824.1111  # H2: 717.4 -> 824.1111
$ws.Cells.Item(2, 9).Value = 695.2105  # I2: 581.6957 -> 695.2105
$ws.Cells.Item(2, 10).Value = 1130.25  # J2: 1163.2858 -> 1130.25
$ws.Cells.Item(2, 11).Value = 695.2105  # K2: 581.6957 -> 695.2105
$ws.Cells.Item(2, 12).Value = 1130.25  # L2: 1163.2858 -> 1130.25
$ws.Cells.Item(2, 13).Value = -582.2105  # M2: -468.6957 -> -582.2105
$ws.Cells.Item(2, 14).Value = -1356.25  # N2: -1389.2858 -> -1356.25

# Row 4: Eyes Bigger than the Plate
$ws.Cells.Item(4, 8).Value = 1133  # H4: 1450 -> 1133
$ws.Cells.Item(4, 9).Value = 949.5  # I4: 0 -> 949.5
$ws.Cells.Item(4, 10).Value = 1500  # J4: 1450 -> 1500
$ws.Cells.Item(4, 11).Value = 949.5  # K4: 0 -> 949.5
$ws.Cells.Item(4, 12).Value = 1500  # L4: 1450 -> 1500
$ws.Cells.Item(4, 13).Value = -833.5  # M4: None -> -833.5
$ws.Cells.Item(4, 14).Value = -1732  # N4: -1682 -> -1732

# Row 10: Bronzed and Burnt
$ws.Cells.Item(10, 8).Value = 400  # H10: 0 -> 400
$ws.Cells.Item(10, 10).Value = 400  # J10: 0 -> 400
$ws.Cells.Item(10, 12).Value = 400  # L10: 0 -> 400
$ws.Cells.Item(10, 14).Value = -740  # N10: None -> -740

# Row 88: The Mast Chance
$ws.Cells.Item(88, 8).Value = 1567.1428  # H88: 1571.8096 -> 1567.1428
$ws.Cells.Item(88, 9).Value = 1036.9  # I88: 1046.7 -> 1036.9
$ws.Cells.Item(88, 11).Value = 1036.9  # K88: 1046.7 -> 1036.9
$ws.Cells.Item(88, 13).Value = -630.9000000000001  # M88: -640.7 -> -630.9000000000001

# Row 91: The Rose and the Riveter (L)
$ws.Cells.Item(91, 8).Value = 1567.1428  # H91: 1571.8096 -> 1567.1428
$ws.Cells.Item(91, 9).Value = 1036.9  # I91: 1046.7 -> 1036.9
$ws.Cells.Item(91, 11).Value = 1036.9  # K91: 1046.7 -> 1036.9
$ws.Cells.Item(91, 13).Value = 367.0999999999999  # M91: 357.3 -> 367.0999999999999

# Row 97: Ore for Me
$ws.Cells.Item(97, 8).Value = 9291.5625  # H97: 10311.786 -> 9291.5625
$ws.Cells.Item(97, 9).Value = 13285.111  # I97: 14683.25 -> 13285.111
$ws.Cells.Item(97, 10).Value = 4157  # J97: 4483.1665 -> 4157
$ws.Cells.Item(97, 11).Value = 13285.111  # K97: 14683.25 -> 13285.111
$ws.Cells.Item(97, 12).Value = 4157  # L97: 4483.1665 -> 4157
$ws.Cells.Item(97, 13).Value = -12789.111  # M97: -14187.25 -> -12789.111
$ws.Cells.Item(97, 14).Value = -5149  # N97: -5475.1665 -> -5149

# Row 102: Smells of Rich Tama-hagane
$ws.Cells.Item(102, 8).Value = 12486.4  # H102: 17196.285 -> 12486.4
$ws.Cells.Item(102, 9).Value = 12486.4  # I102: 19739.334 -> 12486.4
$ws.Cells.Item(102, 10).Value = 0  # J102: 1938 -> 0
$ws.Cells.Item(102, 11).Value = 12486.4  # K102: 19739.334 -> 12486.4
$ws.Cells.Item(102, 12).Value = 0  # L102: 1938 -> 0
$ws.Cells.Item(102, 13).Value = -10864.4  # M102: -18117.334 -> -10864.4
$ws.Cells.Item(102, 14).ClearContents()  # N102: -5182 -> (removed)

# Row 110: Scheduled Maintenance
$ws.Cells.Item(110, 8).Value = 3038.35  # H110: 3079.842 -> 3038.35
$ws.Cells.Item(110, 10).Value = 2145  # J110: 2092.5 -> 2145
$ws.Cells.Item(110, 12).Value = 2145  # L110: 2092.5 -> 2145
$ws.Cells.Item(110, 14).Value = -6235  # N110: -6182.5 -> -6235

# Row 116: No Scope
$ws.Cells.Item(116, 8).Value = 824.1111  # H116: 717.4 -> 824.1111
$ws.Cells.Item(116, 9).Value = 695.2105  # I116: 581.6957 -> 695.2105
$ws.Cells.Item(116, 10).Value = 1130.25  # J116: 1163.2858 -> 1130.25
$ws.Cells.Item(116, 11).Value = 695.2105  # K116: 581.6957 -> 695.2105
$ws.Cells.Item(116, 12).Value = 1130.25  # L116: 1163.2858 -> 1130.25
$ws.Cells.Item(116, 13).Value = 1598.7895  # M116: 1712.3043 -> 1598.7895
$ws.Cells.Item(116, 14).Value = -5718.25  # N116: -5751.2858 -> -5718.25

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Cells.Item(3, 8).Value = 824.1111  # H3: 717.4 -> 824.1111
$ws.Cells.Item(3, 9).Value = 695.2105  # I3: 581.6957 -> 695.2105
$ws.Cells.Item(3, 10).Value = 1130.25  # J3: 1163.2858 -> 1130.25
$ws.Cells.Item(3, 11).Value = 695.2105  # K3: 581.6957 -> 695.2105
$ws.Cells.Item(3, 12).Value = 1130.25  # L3: 1163.2858 -> 1130.25
$ws.Cells.Item(3, 13).Value = -581.2105  # M3: -467.6957 -> -581.2105
$ws.Cells.Item(3, 14).Value = -1358.25  # N3: -1391.2858 -> -1358.25

# Row 20: Smelt and Dealt
$ws.Cells.Item(20, 8).Value = 1497.1818  # H20: 1502.5454 -> 1497.1818
$ws.Cells.Item(20, 9).Value = 1160.625  # I20: 1168 -> 1160.625
$ws.Cells.Item(20, 11).Value = 1160.625  # K20: 1168 -> 1160.625
$ws.Cells.Item(20, 13).Value = -913.625  # M20: -921 -> -913.625

# Row 22: Riveting Run
$ws.Cells.Item(22, 8).Value = 111828.78  # H22: 100666 -> 111828.78
$ws.Cells.Item(22, 9).Value = 500282.5  # I22: 333588.66 -> 500282.5
$ws.Cells.Item(22, 11).Value = 500282.5  # K22: 333588.66 -> 500282.5
$ws.Cells.Item(22, 13).Value = -500109.5  # M22: -333415.66 -> -500109.5

# Row 94: High Steal
$ws.Cells.Item(94, 8).Value = 5306  # H94: 6055.76 -> 5306
$ws.Cells.Item(94, 9).Value = 5969.136  # I94: 6757.25 -> 5969.136
$ws.Cells.Item(94, 10).Value = 2874.5  # J94: 3249.8 -> 2874.5
$ws.Cells.Item(94, 11).Value = 5969.136  # K94: 6757.25 -> 5969.136
$ws.Cells.Item(94, 12).Value = 2874.5  # L94: 3249.8 -> 2874.5
$ws.Cells.Item(94, 13).Value = -5518.136  # M94: -6306.25 -> -5518.136
$ws.Cells.Item(94, 14).Value = -3776.5  # N94: -4151.8 -> -3776.5

# Row 105: Ingot to Wing It
$ws.Cells.Item(105, 8).Value = 3473.805  # H105: 3540.4614 -> 3473.805
$ws.Cells.Item(105, 9).Value = 3400.7632  # I105: 3434.5945 -> 3400.7632
$ws.Cells.Item(105, 10).Value = 4399  # J105: 5499 -> 4399
$ws.Cells.Item(105, 11).Value = 3400.7632  # K105: 3434.5945 -> 3400.7632
$ws.Cells.Item(105, 12).Value = 4399  # L105: 5499 -> 4399
$ws.Cells.Item(105, 13).Value = -1653.7632  # M105: -1687.5945 -> -1653.7632
$ws.Cells.Item(105, 14).Value = -7893  # N105: -8993 -> -7893

# Row 107: The Gold Experience
$ws.Cells.Item(107, 8).Value = 1501.5238  # H107: 1329.4736 -> 1501.5238
$ws.Cells.Item(107, 9).Value = 1238.6842  # I107: 1203.5883 -> 1238.6842
$ws.Cells.Item(107, 10).Value = 3998.5  # J107: 2399.5 -> 3998.5
$ws.Cells.Item(107, 11).Value = 1238.6842  # K107: 1203.5883 -> 1238.6842
$ws.Cells.Item(107, 12).Value = 3998.5  # L107: 2399.5 -> 3998.5
$ws.Cells.Item(107, 13).Value = 681.3158000000001  # M107: 716.4117000000001 -> 681.3158000000001
$ws.Cells.Item(107, 14).Value = -7838.5  # N107: -6239.5 -> -7838.5

$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent
$ws.Cells.Item(7, 8).Value = 90.9375  # H7: 99.64286 -> 90.9375
$ws.Cells.Item(7, 10).Value = 113.625  # J7: 141.5 -> 113.625
$ws.Cells.Item(7, 12).Value = 113.625  # L7: 141.5 -> 113.625
$ws.Cells.Item(7, 14).Value = -339.625  # N7: -367.5 -> -339.625

# Row 31: Wall Not Found
$ws.Cells.Item(31, 8).Value = 6086.4287  # H31: 6319.769 -> 6086.4287
$ws.Cells.Item(31, 9).Value = 3471  # I31: 3755 -> 3471
$ws.Cells.Item(31, 10).Value = 6799.727  # J31: 6786.091 -> 6799.727
$ws.Cells.Item(31, 11).Value = 3471  # K31: 3755 -> 3471
$ws.Cells.Item(31, 12).Value = 6799.727  # L31: 6786.091 -> 6799.727
$ws.Cells.Item(31, 13).Value = -3176  # M31: -3460 -> -3176
$ws.Cells.Item(31, 14).Value = -7389.727  # N31: -7376.091 -> -7389.727

# Row 34: Armoires of the Rich and Famous
$ws.Cells.Item(34, 8).Value = 6086.4287  # H34: 6319.769 -> 6086.4287
$ws.Cells.Item(34, 9).Value = 3471  # I34: 3755 -> 3471
$ws.Cells.Item(34, 10).Value = 6799.727  # J34: 6786.091 -> 6799.727
$ws.Cells.Item(34, 11).Value = 3471  # K34: 3755 -> 3471
$ws.Cells.Item(34, 12).Value = 6799.727  # L34: 6786.091 -> 6799.727
$ws.Cells.Item(34, 13).Value = -3269  # M34: -3553 -> -3269
$ws.Cells.Item(34, 14).Value = -7203.727  # N34: -7190.091 -> -7203.727

# Row 86: Birch, Please
$ws.Cells.Item(86, 8).Value = 20292.5  # H86: 21066.217 -> 20292.5
$ws.Cells.Item(86, 9).Value = 43657.2  # I86: 48230.555 -> 43657.2
$ws.Cells.Item(86, 11).Value = 43657.2  # K86: 48230.555 -> 43657.2
$ws.Cells.Item(86, 13).Value = -42534.2  # M86: -47107.555 -> -42534.2

# Row 89: Built This City on Blocks and Soul (L)
$ws.Cells.Item(89, 8).Value = 20292.5  # H89: 21066.217 -> 20292.5
$ws.Cells.Item(89, 9).Value = 43657.2  # I89: 48230.555 -> 43657.2
$ws.Cells.Item(89, 11).Value = 218286  # K89: 241152.775 -> 218286
$ws.Cells.Item(89, 13).Value = -212670  # M89: -235536.775 -> -212670

# Row 105: Zelkova, My Love
$ws.Cells.Item(105, 8).Value = 570.61536  # H105: 585.1818 -> 570.61536
$ws.Cells.Item(105, 9).Value = 570.61536  # I105: 585.1818 -> 570.61536
$ws.Cells.Item(105, 11).Value = 570.61536  # K105: 585.1818 -> 570.61536
$ws.Cells.Item(105, 13).Value = 1176.38464  # M105: 1161.8182 -> 1176.38464

$ws = $wb.Worksheets.Item("CUL")
# Row 88: Don't Let It Fall Apart
$ws.Cells.Item(88, 8).Value = 22933.066  # H88: 14500 -> 22933.066
$ws.Cells.Item(88, 9).Value = 8998  # I88: 0 -> 8998
$ws.Cells.Item(88, 10).Value = 25076.924  # J88: 14500 -> 25076.924
$ws.Cells.Item(88, 11).Value = 26994  # K88: 0 -> 26994
$ws.Cells.Item(88, 12).Value = 75230.772  # L88: 43500 -> 75230.772
$ws.Cells.Item(88, 13).Value = -26566  # M88: None -> -26566
$ws.Cells.Item(88, 14).Value = -76086.772  # N88: -44356 -> -76086.772

# Row 91: Better Come Back with a Sandwich (L)
$ws.Cells.Item(91, 8).Value = 22933.066  # H91: 14500 -> 22933.066
$ws.Cells.Item(91, 9).Value = 8998  # I91: 0 -> 8998
$ws.Cells.Item(91, 10).Value = 25076.924  # J91: 14500 -> 25076.924
$ws.Cells.Item(91, 11).Value = 26994  # K91: 0 -> 26994
$ws.Cells.Item(91, 12).Value = 75230.772  # L91: 43500 -> 75230.772
$ws.Cells.Item(91, 13).Value = -25512  # M91: None -> -25512
$ws.Cells.Item(91, 14).Value = -78194.772  # N91: -46464 -> -78194.772

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit
$ws.Cells.Item(70, 8).Value = 9000  # H70: 15000 -> 9000
$ws.Cells.Item(70, 9).Value = 8667  # I70: 20001 -> 8667
$ws.Cells.Item(70, 11).Value = 8667  # K70: 20001 -> 8667
$ws.Cells.Item(70, 13).Value = -8397  # M70: -19731 -> -8397

# Row 73: Hulls of Broken Dreams (L)
$ws.Cells.Item(73, 8).Value = 9000  # H73: 15000 -> 9000
$ws.Cells.Item(73, 9).Value = 8667  # I73: 20001 -> 8667
$ws.Cells.Item(73, 11).Value = 8667  # K73: 20001 -> 8667
$ws.Cells.Item(73, 13).Value = -7731  # M73: -19065 -> -7731

# Row 102: Put the Metal to the Peddle
$ws.Cells.Item(102, 8).Value = 2334.6667  # H102: 2451.56 -> 2334.6667
$ws.Cells.Item(102, 9).Value = 1626.5416  # I102: 1695 -> 1626.5416
$ws.Cells.Item(102, 11).Value = 1626.5416  # K102: 1695 -> 1626.5416
$ws.Cells.Item(102, 13).Value = -4.541600000000017  # M102: -73 -> -4.541600000000017

# Row 113: Copious Crystal Cannons
$ws.Cells.Item(113, 8).Value = 4379  # H113: 3199.4285 -> 4379
$ws.Cells.Item(113, 9).Value = 4998  # I113: 2249.5 -> 4998
$ws.Cells.Item(113, 10).Value = 4224.25  # J113: 4466 -> 4224.25
$ws.Cells.Item(113, 11).Value = 4998  # K113: 2249.5 -> 4998
$ws.Cells.Item(113, 12).Value = 4224.25  # L113: 4466 -> 4224.25
$ws.Cells.Item(113, 13).Value = -2828  # M113: -79.5 -> -2828
$ws.Cells.Item(113, 14).Value = -8564.25  # N113: -8806 -> -8564.25

# Row 126: Gold Rush Order
$ws.Cells.Item(126, 8).Value = 5298.5  # H126: 5892 -> 5298.5
$ws.Cells.Item(126, 9).Value = 4476.5  # I126: 4822.6665 -> 4476.5
$ws.Cells.Item(126, 10).Value = 7216.5  # J126: 9100 -> 7216.5
$ws.Cells.Item(126, 11).Value = 13429.5  # K126: 14467.9995 -> 13429.5
$ws.Cells.Item(126, 12).Value = 21649.5  # L126: 27300 -> 21649.5
$ws.Cells.Item(126, 13).Value = -10959.5  # M126: -11997.9995 -> -10959.5
$ws.Cells.Item(126, 14).Value = -26589.5  # N126: -32240 -> -26589.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Cells.Item(22, 8).Value = 50240.81  # H22: 52728.3 -> 50240.81
$ws.Cells.Item(22, 9).Value = 143689  # I22: 143714 -> 143689
$ws.Cells.Item(22, 10).Value = 3516.7144  # J22: 3736 -> 3516.7144
$ws.Cells.Item(22, 11).Value = 143689  # K22: 143714 -> 143689
$ws.Cells.Item(22, 12).Value = 3516.7144  # L22: 3736 -> 3516.7144
$ws.Cells.Item(22, 13).Value = -143394  # M22: -143419 -> -143394
$ws.Cells.Item(22, 14).Value = -4106.7144  # N22: -4326 -> -4106.7144

# Row 27: Fire and Hide
$ws.Cells.Item(27, 8).Value = 50240.81  # H27: 52728.3 -> 50240.81
$ws.Cells.Item(27, 9).Value = 143689  # I27: 143714 -> 143689
$ws.Cells.Item(27, 10).Value = 3516.7144  # J27: 3736 -> 3516.7144
$ws.Cells.Item(27, 11).Value = 143689  # K27: 143714 -> 143689
$ws.Cells.Item(27, 12).Value = 3516.7144  # L27: 3736 -> 3516.7144
$ws.Cells.Item(27, 13).Value = -143582  # M27: -143607 -> -143582
$ws.Cells.Item(27, 14).Value = -3730.7144  # N27: -3950 -> -3730.7144

# Row 68: You Could Say It's a Moving Target
$ws.Cells.Item(68, 8).Value = 3802.6155  # H68: 3923.9285 -> 3802.6155
$ws.Cells.Item(68, 9).Value = 3605  # I68: 3726.7778 -> 3605
$ws.Cells.Item(68, 10).Value = 4118.8  # J68: 4278.8 -> 4118.8
$ws.Cells.Item(68, 11).Value = 3605  # K68: 3726.7778 -> 3605
$ws.Cells.Item(68, 12).Value = 4118.8  # L68: 4278.8 -> 4118.8
$ws.Cells.Item(68, 13).Value = -2856  # M68: -2977.7778 -> -2856
$ws.Cells.Item(68, 14).Value = -5616.8  # N68: -5776.8 -> -5616.8

# Row 71: They Call It Bloody Mary (L)
$ws.Cells.Item(71, 8).Value = 3802.6155  # H71: 3923.9285 -> 3802.6155
$ws.Cells.Item(71, 9).Value = 3605  # I71: 3726.7778 -> 3605
$ws.Cells.Item(71, 10).Value = 4118.8  # J71: 4278.8 -> 4118.8
$ws.Cells.Item(71, 11).Value = 18025  # K71: 18633.889 -> 18025
$ws.Cells.Item(71, 12).Value = 20594  # L71: 21394 -> 20594
$ws.Cells.Item(71, 13).Value = -14281  # M71: -14889.889 -> -14281
$ws.Cells.Item(71, 14).Value = -28082  # N71: -28882 -> -28082
